# Added column for status as of July 4, 2025
# 1. Add a new hidden "DropdownOptions" sheet listing the allowed statuses.
# 2. Add a new header "Status as of July 4, 2025" in AH1 of Sheet1.
# 3. Clear the leftover empty placeholder cells (N2, P2, AE2, AF2).
# 4. Add a list data-validation on AH2 sourced from the DropdownOptions sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Create the hidden DropdownOptions sheet right after Sheet1 ---
$dropdownSheet = $wb.Worksheets.Add($null, $ws)
$dropdownSheet.Name = "DropdownOptions"

$options = @("0% - 10%", "11% - 25%", "26% - 50%", "51% - 75%", "76% - 90%", "91% - 99%", "100%")
for ($i = 0; $i -lt $options.Length; $i++) {
    $cell = $dropdownSheet.Cells.Item($i + 1, 1)
    # Force text storage so values like "100%" aren't coerced into numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $options[$i]
    $cell.Style = "Normal"
}

$dropdownSheet.Visible = $false

# --- 2. New header cell on Sheet1 ---
$ws.Range("AH1").Value = "Status as of July 4, 2025"

# --- 3. Drop the stray empty cells left in row 2 ---
$ws.Range("N2").ClearContents()
$ws.Range("P2").ClearContents()
$ws.Range("AE2").ClearContents()
$ws.Range("AF2").ClearContents()

# --- 4. Data validation (dropdown list) on AH2 ---
$validation = $ws.Range("AH2").Validation
$validation.Add(3, 1, 1, "DropdownOptions!`$A`$1:`$A`$7")
$validation.IgnoreBlank = $true
$validation.InCellDropdown = $true
$validation.ShowInput = $false
$validation.ShowError = $false

# Leave the workbook focused back on the original sheet.
$ws.Activate()
